$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = -2
    4  = 1
    6  = 3
    7  = -2
    8  = 2
    9  = -3
    10 = 2
    11 = 8
    13 = 4
    14 = 3
    15 = 8
    17 = 3
    18 = -1
    19 = 2
    20 = 3
    21 = 2
    22 = -1
    23 = -4
    25 = -1
    26 = -3
    27 = 5
    28 = 3
    29 = 1
    31 = 6
    32 = 4
    33 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 6).Value = $changes[$row]
}
